$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend the "USE CASE DESCRIPTIONS" table (rows 3-10) with an extra
#    "Check 1" column (E) for rows 11-17, matching the pattern already used
#    by rows 3-10 (column E, fillId style used on E3:E10).
# ---------------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("E11:E17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Add a new "CLASS DIAGRAM" heading in row 19 (new section, plain text,
#    no special styling - matches B3's un-styled heading cell).
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "CLASS DIAGRAM"

# ---------------------------------------------------------------------------
# 3. Mark a "highlight" cell next to ACTIVITY DIAGRAMS (row 21, column D)
#    using the same green fill already used for column D (copied from D3)
#    plus a red font color - Excel will create a new font + cell style for
#    this combination automatically.
# ---------------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D21").Font.Color = 255

# ---------------------------------------------------------------------------
# 4. Add a new "SEQUENCE DIAGRAM" section (rows 23-29), re-using the same
#    use-case names/style already used for the ACTIVITY DIAGRAMS list
#    (rows 11-17, column C).
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "SEQUENCE DIAGRAM"

$useCases = @("Register", "AdministerArticles", "BrowseArticles", "GetRecommendations", "InteractWithArticles", "ManageProfile", "ManageUserAccounts")
$row = 23
foreach ($name in $useCases) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $name
    $cell.HorizontalAlignment = -4131
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 5. Misc sheet view / print tweaks.
# ---------------------------------------------------------------------------
$null = $ws.Range("E17").Select()
$ws.PageSetup.Orientation = 1
